$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits right after "Test2" in the second
# paragraph. Remove it from there -- it will be re-created later at the end
# of the document.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Append a new paragraph containing "Hello Everyone" plus a further new
# (empty) paragraph that holds the "_GoBack" bookmark, after the existing
# trailing empty paragraph at the end of the document body.
$endRange = $d.Range($d.Content.End, $d.Content.End)

$xmlFragment = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p/>
          <w:p><w:r><w:t>Hello Everyone</w:t></w:r></w:p>
          <w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$endRange.InsertXML($xmlFragment)
